$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: period headers (rolling window - drop oldest, append newest)
$ws.Range("D8").Value = "3 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "6 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "9 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "3 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "6 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "9 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "12 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "3 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "6 ماهه منتهی به 1401/12"

# Row 9: publish dates
$ws.Range("D9").Value = "1400-10-29 (2)"
$ws.Range("E9").Value = "1401-03-21 (4)"
$ws.Range("F9").Value = "1401-05-19 (3)"
$ws.Range("G9").Value = "1401-10-05 (9)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-01-30 (3)"
$ws.Range("J9").Value = "1401-05-19 (2)"
$ws.Range("K9").Value = "1402-01-30 (4)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-01-30"

# Row 11
$ws.Range("D11").Value = 46810537
$ws.Range("E11").Value = 88806190
$ws.Range("F11").Value = 136740382
$ws.Range("G11").Value = 197398681
$ws.Range("H11").Value = 94387382
$ws.Range("I11").Value = 161373844
$ws.Range("J11").Value = 289433821
$ws.Range("K11").Value = 431582823
$ws.Range("L11").Value = 110842812
$ws.Range("M11").Value = 181186833

# Row 12
$ws.Range("D12").Value = -14466354
$ws.Range("E12").Value = -30126728
$ws.Range("F12").Value = -44214185
$ws.Range("G12").Value = -70710577
$ws.Range("H12").Value = -37794904
$ws.Range("I12").Value = -68923090
$ws.Range("J12").Value = -110912374
$ws.Range("K12").Value = -159682555
$ws.Range("L12").Value = -34318064
$ws.Range("M12").Value = -65203331

# Row 13
$ws.Range("D13").Value = 32344183
$ws.Range("E13").Value = 58679462
$ws.Range("F13").Value = 92526197
$ws.Range("G13").Value = 126688104
$ws.Range("H13").Value = 56592478
$ws.Range("I13").Value = 92450754
$ws.Range("J13").Value = 178521447
$ws.Range("K13").Value = 271900268
$ws.Range("L13").Value = 76524748
$ws.Range("M13").Value = 115983502

# Row 14
$ws.Range("D14").Value = -3831223
$ws.Range("E14").Value = -7656550
$ws.Range("F14").Value = -10794753
$ws.Range("G14").Value = -9591839
$ws.Range("H14").Value = -2265379
$ws.Range("I14").Value = -3606877
$ws.Range("J14").Value = -5279200
$ws.Range("K14").Value = -38847143
$ws.Range("L14").Value = -3191115
$ws.Range("M14").Value = -46857925

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

# Row 16
$ws.Range("D16").Value = 273315
$ws.Range("E16").Value = 3938522
$ws.Range("F16").Value = 1565354
$ws.Range("G16").Value = 2372391
$ws.Range("H16").Value = -169830
$ws.Range("I16").Value = 1467698
$ws.Range("J16").Value = 1026848
$ws.Range("K16").Value = 11254375
$ws.Range("L16").Value = 1836195
$ws.Range("M16").Value = 23315668

# Row 17
$ws.Range("D17").Value = 28786275
$ws.Range("E17").Value = 54961434
$ws.Range("F17").Value = 83296798
$ws.Range("G17").Value = 119468656
$ws.Range("H17").Value = 54157269
$ws.Range("I17").Value = 90311575
$ws.Range("J17").Value = 174269095
$ws.Range("K17").Value = 244307500
$ws.Range("L17").Value = 75169828
$ws.Range("M17").Value = 92441245

# Row 18
$ws.Range("D18").Value = -77409
$ws.Range("E18").Value = -249766
$ws.Range("F18").Value = -316159
$ws.Range("G18").Value = -474078
$ws.Range("H18").Value = -258982
$ws.Range("I18").Value = -438268
$ws.Range("J18").Value = -588232
$ws.Range("K18").Value = -1065627
$ws.Range("L18").Value = -484883
$ws.Range("M18").Value = -954828

# Row 19
$ws.Range("D19").Value = 115810
$ws.Range("E19").Value = 1147656
$ws.Range("F19").Value = 2348964
$ws.Range("G19").Value = 4525226
$ws.Range("H19").Value = 2410362
$ws.Range("I19").Value = 6280618
$ws.Range("J19").Value = 9101554
$ws.Range("K19").Value = 12770725
$ws.Range("L19").Value = 3330336
$ws.Range("M19").Value = 15586152

# Row 20
$ws.Range("D20").Value = 28824676
$ws.Range("E20").Value = 55859324
$ws.Range("F20").Value = 85329603
$ws.Range("G20").Value = 123519804
$ws.Range("H20").Value = 56308649
$ws.Range("I20").Value = 96153925
$ws.Range("J20").Value = 182782417
$ws.Range("K20").Value = 256012598
$ws.Range("L20").Value = 78015281
$ws.Range("M20").Value = 107072569

# Row 21
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = -26316
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 4597514
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = -16130
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = -184701
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -41198

# Row 22
$ws.Range("D22").Value = 28824676
$ws.Range("E22").Value = 55833008
$ws.Range("F22").Value = 85329603
$ws.Range("G22").Value = 128117318
$ws.Range("H22").Value = 56308649
$ws.Range("I22").Value = 96137795
$ws.Range("J22").Value = 182782417
$ws.Range("K22").Value = 255827897
$ws.Range("L22").Value = 78015281
$ws.Range("M22").Value = 107031371

# Row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("D24").Value = 28824676
$ws.Range("E24").Value = 55833008
$ws.Range("F24").Value = 85329603
$ws.Range("G24").Value = 128117318
$ws.Range("H24").Value = 56308649
$ws.Range("I24").Value = 96137795
$ws.Range("J24").Value = 182782417
$ws.Range("K24").Value = 255827897
$ws.Range("L24").Value = 78015281
$ws.Range("M24").Value = 107031371

# Row 25
$ws.Range("D25").Value = 4804
$ws.Range("E25").Value = 9306
$ws.Range("F25").Value = 14222
$ws.Range("G25").Value = 21353
$ws.Range("H25").Value = 9385
$ws.Range("I25").Value = 16023
$ws.Range("J25").Value = 30464
$ws.Range("K25").Value = 42638
$ws.Range("L25").Value = 13003
$ws.Range("M25").Value = 17839

# Row 26
$ws.Range("D26").Value = 6000000
$ws.Range("E26").Value = 6000000
$ws.Range("F26").Value = 6000000
$ws.Range("G26").Value = 6000000
$ws.Range("H26").Value = 6000000
$ws.Range("I26").Value = 6000000
$ws.Range("J26").Value = 6000000
$ws.Range("K26").Value = 6000000
$ws.Range("L26").Value = 6000000
$ws.Range("M26").Value = 6000000

# Row 27
$ws.Range("D27").Value = 4804
$ws.Range("E27").Value = 9306
$ws.Range("F27").Value = 14222
$ws.Range("G27").Value = 21353
$ws.Range("H27").Value = 9385
$ws.Range("I27").Value = 16023
$ws.Range("J27").Value = 30464
$ws.Range("K27").Value = 42638
$ws.Range("L27").Value = 13003
$ws.Range("M27").Value = 17839

# Column width adjustments (shift left pattern)
$ws.Range("D:D").ColumnWidth = 27.1666666667
$ws.Range("G:G").ColumnWidth = 28.1666666667
$ws.Range("H:H").ColumnWidth = 27.1666666667
$ws.Range("K:K").ColumnWidth = 28.1666666667
$ws.Range("L:L").ColumnWidth = 27.1666666667